$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 22950771235.97675
$ws.Range("D3").Value = 24375759418.22945
$ws.Range("E3").Value = 89718.73
$ws.Range("F3").Value = 89675
